# Weekly CompStat update for the 104th Precinct.
# Report period moves forward one week:
#   "Number 33" -> "Number 34"
#   8/12/2024-8/18/2024 -> 8/19/2024-8/25/2024
# plus the refreshed crime-statistics figures for the new week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Title / header text
# ---------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  34"
$ws.Range("C9").Value = "Report Covering the Week  8/19/2024  Through  8/25/2024"

# ---------------------------------------------------------------
# Helper donor cells used to re-create the correct cell style when a
# cell needs to switch between a numeric value and the "N/A" shared
# text (style 14, shared-string "0" or "***.*").
#   NA0 donor (text "0")     -> C23
#   NA1 donor (text "***.*") -> E23
#   Numeric "count" style donor (s=15) -> I31
#   Numeric "pct"   style donor (s=16) -> L31
# ---------------------------------------------------------------
$naZero = $ws.Range("C23")
$naDash = $ws.Range("E23")
$numStyle = $ws.Range("I31")
$pctStyle = $ws.Range("L31")

function Set-NA0([string]$addr) {
    $naZero.Copy($ws.Range($addr))
}
function Set-NA1([string]$addr) {
    $naDash.Copy($ws.Range($addr))
}
function Set-Num([string]$addr, $value) {
    $numStyle.Copy($ws.Range($addr))
    $ws.Range($addr).Value = $value
}
function Set-Pct([string]$addr, $value) {
    $pctStyle.Copy($ws.Range($addr))
    $ws.Range($addr).Value = $value
}

# ---------------------------------------------------------------
# Row 14 - Murder : Week-to-date columns (G/H) become N/A
# ---------------------------------------------------------------
Set-NA0 "G14"
Set-NA1 "H14"

# ---------------------------------------------------------------
# Row 15 - Rape : 28-day prior (D/E) become numeric, WTD updates
# ---------------------------------------------------------------
Set-Num "D15" 2
Set-Pct "E15" -100
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 14
$ws.Range("K15").Value = 14.285714285714

# ---------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = -21.052631578947
$ws.Range("I16").Value = 137
$ws.Range("J16").Value = 142
$ws.Range("K16").Value = -3.521126760563
$ws.Range("L16").Value = 6.201550387596
$ws.Range("M16").Value = -20.348837209302
$ws.Range("N16").Value = -77.577741407528

# ---------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 61.538461538461
$ws.Range("I17").Value = 197
$ws.Range("J17").Value = 181
$ws.Range("K17").Value = 8.839779005524
$ws.Range("L17").Value = 2.072538860103
$ws.Range("M17").Value = 31.333333333333
$ws.Range("N17").Value = -2.955665024630

# ---------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 30.769230769230
$ws.Range("I18").Value = 150
$ws.Range("J18").Value = 123
$ws.Range("K18").Value = 21.951219512195
$ws.Range("L18").Value = -17.127071823204
$ws.Range("M18").Value = -50.331125827814
$ws.Range("N18").Value = -88.986784140969

# ---------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = -23.728813559322
$ws.Range("I19").Value = 436
$ws.Range("J19").Value = 445
$ws.Range("K19").Value = -2.022471910112
$ws.Range("L19").Value = 5.060240963855
$ws.Range("M19").Value = 60.885608856088
$ws.Range("N19").Value = 1.160092807424

# ---------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 14
$ws.Range("E20").Value = -57.142857142857
$ws.Range("F20").Value = 42
$ws.Range("G20").Value = 38
$ws.Range("H20").Value = 10.526315789473
$ws.Range("I20").Value = 273
$ws.Range("J20").Value = 236
$ws.Range("K20").Value = 15.677966101694
$ws.Range("L20").Value = 55.113636363636
$ws.Range("M20").Value = 10.080645161290
$ws.Range("N20").Value = -88.471283783783

# ---------------------------------------------------------------
# Row 21 - TOTAL (D21 stays 37, unchanged)
# ---------------------------------------------------------------
$ws.Range("C21").Value = 31
$ws.Range("E21").Value = -16.216216216216
$ws.Range("F21").Value = 140
$ws.Range("G21").Value = 144
$ws.Range("H21").Value = -2.777777777777
$ws.Range("I21").Value = 1210
$ws.Range("J21").Value = 1144
$ws.Range("K21").Value = 5.769230769230
$ws.Range("L21").Value = 9.107303877367
$ws.Range("M21").Value = 4.671280276816
$ws.Range("N21").Value = -75.843481732880

# ---------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------
$ws.Range("G22").Value = 1

# ---------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = -37.142857142857
$ws.Range("F24").Value = 107
$ws.Range("G24").Value = 114
$ws.Range("H24").Value = -6.140350877192
$ws.Range("I24").Value = 938
$ws.Range("J24").Value = 857
$ws.Range("K24").Value = 9.451575262543
$ws.Range("L24").Value = 0.213675213675
$ws.Range("M24").Value = 34.964028776978

# ---------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = -23.076923076923
$ws.Range("F25").Value = 54
$ws.Range("G25").Value = 46
$ws.Range("H25").Value = 17.391304347826
$ws.Range("I25").Value = 362
$ws.Range("J25").Value = 297
$ws.Range("K25").Value = 21.885521885521
$ws.Range("L25").Value = -5.729166666666

# ---------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -20
$ws.Range("F26").Value = 43
$ws.Range("G26").Value = 27
$ws.Range("H26").Value = 59.259259259259
$ws.Range("I26").Value = 390
$ws.Range("J26").Value = 310
$ws.Range("K26").Value = 25.806451612903
$ws.Range("L26").Value = 13.702623906705
$ws.Range("M26").Value = -18.067226890756

# ---------------------------------------------------------------
# Row 27 - UCR Rape* : D/E become numeric, F becomes N/A
# ---------------------------------------------------------------
Set-Num "D27" 2
Set-Pct "E27" -100
Set-NA0 "F27"
$ws.Range("H27").Value = -100
$ws.Range("J27").Value = 23
$ws.Range("K27").Value = -4.347826086956
$ws.Range("L27").Value = 0

# ---------------------------------------------------------------
# Row 28 - Other Sex Crimes : C becomes numeric
# ---------------------------------------------------------------
Set-Num "C28" 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 22
$ws.Range("J28").Value = 45
$ws.Range("K28").Value = -51.111111111111
$ws.Range("L28").Value = -56.862745098039

# ---------------------------------------------------------------
# Row 29 - Shooting Vic. : G/H become N/A
# ---------------------------------------------------------------
Set-NA0 "G29"
Set-NA1 "H29"
$ws.Range("N29").Value = -90.909090909090

# ---------------------------------------------------------------
# Row 30 - Shooting Inc. : G/H become N/A
# ---------------------------------------------------------------
Set-NA0 "G30"
Set-NA1 "H30"
$ws.Range("N30").Value = -89.473684210526

# ---------------------------------------------------------------
# Row 33 - Traffic Fatalities : D/E become N/A
# ---------------------------------------------------------------
Set-NA0 "D33"
Set-NA1 "E33"
$ws.Range("G33").Value = 3
